$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Range("H33").Value = 454.33334
$ws.Range("I33").Value = 479.58334
$ws.Range("K33").Value = 479.58334
$ws.Range("M33").Value = -250.58334
# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1939.7091
$ws.Range("I138").Value = 2054.3635
$ws.Range("J138").Value = 1911.0454
$ws.Range("K138").Value = 6163.0905
$ws.Range("L138").Value = 5733.1362
$ws.Range("M138").Value = -1023.0905
$ws.Range("N138").Value = -16013.1362

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 2864.324
$ws.Range("I32").Value = 2685.7969
$ws.Range("K32").Value = 2685.7969
$ws.Range("M32").Value = -2398.7969
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 1460.0588
$ws.Range("I61").Value = 1080.9
$ws.Range("K61").Value = 1080.9
$ws.Range("M61").Value = -868.9000000000001
# Row 62: Hauberk and No Play | Mythrite Hauberk of Maiming
$ws.Range("H62").Value = 65000
$ws.Range("J62").Value = 65000
$ws.Range("L62").Value = 65000
$ws.Range("N62").Value = -66248
# Row 65: Knights without Armor (L) | Mythrite Hauberk of Maiming
$ws.Range("H65").Value = 65000
$ws.Range("J65").Value = 65000
$ws.Range("L65").Value = 195000
$ws.Range("N65").Value = -201240
# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 1047.5834
$ws.Range("I110").Value = 823.8333
$ws.Range("K110").Value = 823.8333
$ws.Range("M110").Value = 1221.1667
# Row 121: Shield to Shield | Dwarven Mythril Shield
$ws.Range("H121").Value = 36451
$ws.Range("J121").Value = 36451
$ws.Range("L121").Value = 36451
$ws.Range("N121").Value = -39945
# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1460.0588
$ws.Range("I136").Value = 1080.9
$ws.Range("K136").Value = 3242.7
$ws.Range("M136").Value = -692.7000000000003

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 6533.773
$ws.Range("I134").Value = 1152.3889
$ws.Range("J134").Value = 30750
$ws.Range("K134").Value = 3457.1667
$ws.Range("L134").Value = 92250
$ws.Range("M134").Value = -922.1666999999998
$ws.Range("N134").Value = -97320

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 17: Say It with Spears | Feathered Harpoon
$ws.Range("H17").Value = 3300.5
$ws.Range("I17").Value = 2600
$ws.Range("J17").Value = 4001
$ws.Range("K17").Value = 2600
$ws.Range("L17").Value = 4001
$ws.Range("M17").Value = -2426
$ws.Range("N17").Value = -4349
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2029.4348
$ws.Range("I31").Value = 1751.3077
$ws.Range("K31").Value = 1751.3077
$ws.Range("M31").Value = -1456.3077
# Row 33: Tools for the Tools | Silver Battle Fork
$ws.Range("H33").Value = 1049.5
$ws.Range("I33").Value = 1049.5
$ws.Range("K33").Value = 1049.5
$ws.Range("M33").Value = -670.5
# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2029.4348
$ws.Range("I34").Value = 1751.3077
$ws.Range("K34").Value = 1751.3077
$ws.Range("M34").Value = -1549.3077
# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 8270.833000000001
$ws.Range("I132").Value = 12116.3
$ws.Range("J132").Value = 3464
$ws.Range("K132").Value = 36348.89999999999
$ws.Range("L132").Value = 10392
$ws.Range("M132").Value = -33818.89999999999
$ws.Range("N132").Value = -15452
# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 1981.1482
$ws.Range("I134").Value = 1953.6666
$ws.Range("K134").Value = 5860.9998
$ws.Range("M134").Value = -3325.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 116: On a Full Stomach | Sausage Links
$ws.Range("H116").Value = 2899.8
$ws.Range("I116").Value = 2249.5
$ws.Range("J116").Value = 3333.3333
$ws.Range("K116").Value = 6748.5
$ws.Range("L116").Value = 9999.999899999999
$ws.Range("M116").Value = -3306.5
$ws.Range("N116").Value = -16883.9999
# Row 136: Simple Is Hardest | Spaghetti al Olio e Peperoncino
$ws.Range("H136").Value = 1636.5555
$ws.Range("I136").Value = 1406
$ws.Range("J136").Value = 1924.75
$ws.Range("K136").Value = 4218
$ws.Range("L136").Value = 5774.25
$ws.Range("M136").Value = 882
$ws.Range("N136").Value = -15974.25
# Row 139: Najoothie | Wild Banana Blend
$ws.Range("H139").Value = 1573.1794
$ws.Range("I139").Value = 1625.7916
$ws.Range("K139").Value = 4877.3748
$ws.Range("M139").Value = 262.6252000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 40912708
$ws.Range("I70").Value = 31253774
$ws.Range("K70").Value = 31253774
$ws.Range("M70").Value = -31253504
# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 40912708
$ws.Range("I73").Value = 31253774
$ws.Range("K73").Value = 31253774
$ws.Range("M73").Value = -31252838
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 3529.5
$ws.Range("J80").Value = 3466.6667
$ws.Range("L80").Value = 3466.6667
$ws.Range("N80").Value = -5462.6667
# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 3529.5
$ws.Range("J83").Value = 3466.6667
$ws.Range("L83").Value = 17333.3335
$ws.Range("N83").Value = -27317.3335
# Row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Range("H107").Value = 776.56525
$ws.Range("J107").Value = 650.8333
$ws.Range("L107").Value = 650.8333
$ws.Range("N107").Value = -4490.8333
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 1355.4375
$ws.Range("I122").Value = 1461.125
$ws.Range("J122").Value = 1038.375
$ws.Range("K122").Value = 4383.375
$ws.Range("L122").Value = 3115.125
$ws.Range("M122").Value = -1933.375
$ws.Range("N122").Value = -8015.125
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 2016.25
$ws.Range("I126").Value = 1670.0769
$ws.Range("K126").Value = 5010.2307
$ws.Range("M126").Value = -2540.2307
# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 1781.4054
$ws.Range("I132").Value = 1668.6
$ws.Range("J132").Value = 2016.4166
$ws.Range("K132").Value = 5005.799999999999
$ws.Range("L132").Value = 6049.2498
$ws.Range("M132").Value = -2475.799999999999
$ws.Range("N132").Value = -11109.2498

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 5458.8423
$ws.Range("I46").Value = 816
$ws.Range("J46").Value = 8835.454
$ws.Range("K46").Value = 816
$ws.Range("L46").Value = 8835.454
$ws.Range("M46").Value = -628
$ws.Range("N46").Value = -9211.454

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 66667910
$ws.Range("I126").Value = 71429870
$ws.Range("J126").Value = 500
$ws.Range("K126").Value = 214289610
$ws.Range("L126").Value = 1500
$ws.Range("M126").Value = -214287140
$ws.Range("N126").Value = -6440
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1510.4348
$ws.Range("I132").Value = 1027.8182
$ws.Range("J132").Value = 2735.5386
$ws.Range("K132").Value = 3083.4546
$ws.Range("L132").Value = 8206.6158
$ws.Range("M132").Value = -553.4546
$ws.Range("N132").Value = -13266.6158
# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 399.5484
$ws.Range("I136").Value = 303.28
$ws.Range("J136").Value = 800.6667
$ws.Range("K136").Value = 909.8399999999999
$ws.Range("L136").Value = 2402.0001
$ws.Range("M136").Value = 1640.16
$ws.Range("N136").Value = -7502.0001
